# Update "Training Dashboard" sheet rows 3-22:
#   - Column H ("PERIOD TO EXPIRE"): decrement by 1
#   - Column I ("LAST UPDATE"): change date text from 03-Nov-2025 to 04-Nov-2025

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

for ($row = 3; $row -le 22; $row++) {
    $hCell = $ws.Cells.Item($row, 8)   # column H
    $hCell.Value2 = $hCell.Value2 - 1

    $iCell = $ws.Cells.Item($row, 9)   # column I
    # Leading apostrophe forces the date-looking text to stay a literal
    # string instead of being parsed into a date serial number.
    $iCell.Value2 = "'04-Nov-2025"
}
